$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing odds on row 96 (Hatayspor vs Adana Demirspor) ---
$ws.Range("F96").Value = 2.7
$ws.Range("G96").Value = 3.6
$ws.Range("H96").Value = 2.42
$ws.Range("K96").Value = 1.7

# --- Update existing odds on row 122 (Antalyaspor vs Kasimpasa) ---
$ws.Range("F122").Value = 1.78
$ws.Range("G122").Value = 3.6
$ws.Range("H122").Value = 4.55
$ws.Range("K122").Value = 1.79

# --- Update existing odds on row 123 (Ankaragucu vs Besiktas) ---
$ws.Range("F123").Value = 3.5
$ws.Range("G123").Value = 3.45
$ws.Range("H123").Value = 2.05
$ws.Range("K123").Value = 1.65

# --- Append new row 147 (Seattle Sounders vs Houston Dynamo) ---
$ws.Range("A147").Value = "USA MLS"
$ws.Range("B147").Value = 44808.91666666666
$ws.Range("B147").NumberFormat = $ws.Range("B146").NumberFormat
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = "Seattle Sounders"
$ws.Range("E147").Value = "Houston Dynamo"
$ws.Range("F147").Value = 1.81
$ws.Range("G147").Value = 3.93
$ws.Range("H147").Value = 3.82
$ws.Range("I147").Value = 1.02
$ws.Range("J147").Value = 1.2
$ws.Range("K147").Value = 1.6
$ws.Range("L147").Value = 1.62
$ws.Range("M147").Value = 2.15
$ws.Range("N147").Value = 1.3

# --- Append new row 148 (SJ Earthquakes vs Vancouver Whitecaps) ---
$ws.Range("A148").Value = "USA MLS"
$ws.Range("B148").Value = 44808.9375
$ws.Range("B148").NumberFormat = $ws.Range("B146").NumberFormat
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = "SJ Earthquakes"
$ws.Range("E148").Value = "Vancouver Whitecaps"
$ws.Range("F148").Value = 1.96
$ws.Range("G148").Value = 3.65
$ws.Range("H148").Value = 3.75
$ws.Range("I148").Value = 1.02
$ws.Range("J148").Value = 1.18
$ws.Range("K148").Value = 1.56
$ws.Range("L148").Value = 1.57
$ws.Range("M148").Value = 2.25
$ws.Range("N148").Value = 1.32
